# Add PF/1.0.5 to meta-sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.5"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row uses the default/Normal style (no inherited column formatting),
# matching the rest of the sheet's "Normal" cell style.
$ws.Range("A3:D3").Style = "Normal"
